$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the '01.03.2024' sheet's input figures (dt_mitgcm, date_start,
#    date_end, ID_restart_file). Dependent formulas (C3, G8, G9, G11, C12)
#    recalc automatically.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("01.03.2024")
$ws4.Range("G3").Value = 32
$ws4.Range("G4").Value = 45512
$ws4.Range("G5").Value = 45658
$ws4.Range("C6").Value = 432000
$ws4.Range("C2:C3").Select()

# ---------------------------------------------------------------------
# 2) Duplicate the '01.12.2023' sheet to create the new '01.03.2025' tab,
#    inserting it right before '01.12.2023'.
# ---------------------------------------------------------------------
$ws12old = $wb.Worksheets.Item("01.12.2023")
$ws12old.Copy($ws12old)
$wsnew = $wb.ActiveSheet
$wsnew.Name = "01.03.2025"

# Restore the plain selection on the original '01.12.2023' sheet (re-fetch
# by name since the old object reference now resolves to the new copy).
$ws12 = $wb.Worksheets.Item("01.12.2023")
$ws12.Range("A1:G12").Select()

# ---------------------------------------------------------------------
# 3) Fill in the new sheet's own figures, then leave it the active tab.
# ---------------------------------------------------------------------
$wsnew = $wb.Worksheets.Item("01.03.2025")
$wsnew.Activate()
$wsnew.Range("I5").Clear()
$wsnew.Range("G2").Value = 45778
$wsnew.Range("C2").Formula = "=G2"
$wsnew.Range("G3").Value = 32
$wsnew.Range("C3").Formula = "=G3"
$wsnew.Range("C6").Value = 396900
$wsnew.Range("K6").Select()
